$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 / 7: swap Starting/Ending SoC (%) values ---
$ws.Range("B6").Value = 99
$ws.Range("B7").Value = 10

# --- Row 8: distance label ---
$ws.Range("A8").Value = "Total distance covered (km)"

# --- Row 9: WH/KM label ---
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"

# --- Row 10: Total SOC consumed ---
$ws.Range("A10").Value = "Total SOC consumed(%)"
$ws.Range("B10").Value = 89

# --- Row 12: Peak Power ---
$ws.Range("A12").Value = "Peak Power(kW)"

# --- Row 13: Average Power ---
$ws.Range("A13").Value = "Average Power(kW)"

# --- Row 14: Total Energy Regenerated ---
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"

# --- Row 15: Regenerative Effectiveness ---
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 6.855366834693272

# --- Row 16 / 17: swap Lowest/Highest Cell Voltage (label + value) ---
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.376
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.106

# --- Row 18: Difference in Cell Voltage ---
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"

# --- Row 19: Minimum Temperature ---
$ws.Range("A19").Value = "Minimum Temperature(C)"

# --- Row 20: Maximum Temperature ---
$ws.Range("A20").Value = "Maximum Temperature(C)"

# --- Row 21: Difference in Temperature ---
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 10

# --- Row 22: Maximum Fet Temperature ---
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"

# --- Row 23: Maximum Afe Temperature ---
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"

# --- Row 24: Maximum PCB Temperature ---
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"

# --- Row 25: Maximum MCU Temperature ---
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"

# --- Row 26: Maximum Motor Temperature ---
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"

# --- Row 27: Abnormal Motor Temperature Detected ---
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

# --- Row 28 / 29: swap lowest/highest cell temp labels ---
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("A29").Value = "lowest cell temp(C)"

# --- Row 30: Difference between Highest and Lowest Cell Temperature at 100% SOC ---
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# --- Row 31: was "Maximum BMS Temperature in C" -> "Battery Voltage(V)" ---
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 55

# --- Row 32: was "Battery Voltage" -> "Total energy charged(kWh)" ---
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.948992665277778

# --- Row 33: was "Total energy charged in kWh" -> "Electricity consumption units(kW)" ---
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.00000007890786349891408

# --- Row 34: was "Electricity consumption units in kW" -> "Idling time percentage" ---
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 5.18963018963019

# --- Row 35: was "Idling time percentage" -> "Time spent in 0-10 km/h" ---
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 10.28948528948529

# --- Row 36: was "Time spent in 0-10 km/h" -> "Time spent in 10-20 km/h" ---
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 7.397782397782398

# --- Row 37: was "Time spent in 10-20 km/h" -> "Time spent in 20-30 km/h" ---
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 11.45498645498645

# --- Row 38: was "Time spent in 20-30 km/h" -> "Time spent in 30-40 km/h" ---
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 44.86549486549487

# --- Row 39: was "Time spent in 30-40 km/h" -> "Time spent in 40-50 km/h" ---
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 18.42909342909343

# --- Row 40: was "Time spent in 40-50 km/h" -> "Time spent in 50-60 km/h" ---
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 0.2488502488502489

# --- Row 41: was "Time spent in 50-60 km/h" -> "Time spent in 60-70 km/h" ---
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 0.3858753858753859

# --- Row 42: was "Time spent in 60-70 km/h" -> "Time spent in 70-80 km/h" ---
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 0.3087003087003087

# --- Row 43 (new row) ---
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
